$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.474.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.49%  "

$ws.Range("D3").Value = "'3.066.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.80%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'551.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.07%  "

$ws.Range("D6").Value = "'142.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.37%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'3.062.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.91%  "

$ws.Range("D9").Value = "'0.503"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.60%  "

$ws.Range("D10").Value = "'6.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.53%  "

$ws.Range("E11").Value = "  +2.72%  "

$ws.Range("D12").Value = "'0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.09%  "

$ws.Range("E13").Value = "  +2.81%  "

$ws.Range("D14").Value = "'34.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.66%  "

$ws.Range("D15").Value = "'3.565.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.89%  "

$ws.Range("D16").Value = "'63.427.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.38%  "

$ws.Range("D17").Value = "'3.067.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.77%  "

$ws.Range("E19").Value = "  +3.03%  "

$ws.Range("D20").Value = "'483.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.10%  "

$ws.Range("D21").Value = "'13.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.60%  "

$ws.Range("D22").Value = "'0.676"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.20%  "

$ws.Range("D23").Value = "'7.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.98%  "

$ws.Range("D24").Value = "'80.74"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'12.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.63%  "

$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("E27").Value = "  +4.30%  "

$ws.Range("D28").Value = "'7.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.13%  "

$ws.Range("D29").Value = "'2.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.77%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").Value = "  +2.71%  "

$ws.Range("E32").Value = "  +1.58%  "

$ws.Range("E33").Value = "  +8.04%  "

$ws.Range("E34").Value = "  +5.50%  "

$ws.Range("D35").Value = "'55.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("D36").Value = "'5.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.62%  "

$ws.Range("D37").Value = "'465.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.29%  "

$ws.Range("D38").Value = "'0.0822"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.26%  "

$ws.Range("D39").Value = "'0.0395"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.99%  "

$ws.Range("D40").Value = "'3.007.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.40%  "

$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("E42").Value = "  +1.95%  "

$ws.Range("E43").Value = "  +5.93%  "

$ws.Range("D44").Value = "'27.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.92%  "

$ws.Range("E45").Value = "  +5.84%  "

$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("E47").Value = "  +3.07%  "

$ws.Range("E48").Value = "  +2.75%  "

$ws.Range("D49").Value = "'116.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "

$ws.Range("E50").Value = "  +3.32%  "

$ws.Range("E51").Value = "  +4.67%  "
